$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: columns I, J, K, L
$ws.Range("I2").Value = -0.7223882864291703
$ws.Range("J2").Value = 0.2865576572302282
$ws.Range("K2").Value = 0.1359011226912011
$ws.Range("L2").Value = 2.611494526868577

# Row 19: columns I, J, K, L
$ws.Range("I19").Value = -0.9527902356820649
$ws.Range("J19").Value = 0.3475976934985399
$ws.Range("K19").Value = 0.4209129885157118
$ws.Range("L19").Value = 2.411671512147268
